# Generate Report for Handback
# Adds a new handback entry (563683f9-0bfd-4622-870a-0da819665c33) as row 4
# on the "Overview", "zh-cn" and "de-de" worksheets, mirroring the existing
# rows 2/3 pattern (d82cf695-..., 298923d6-...).

$wb = $excel.ActiveWorkbook

$uuid  = "563683f9-0bfd-4622-870a-0da819665c33"
$token = "407048bd854fcf0d97e6a44f94b8fc8b83b13639"
$md    = "$uuid.md"

# ---------------------------------------------------------------------
# Sheet "Overview": columns File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/407048bd854fcf0d97e6a44f94b8fc8b83b13639/e2e/$md", "", "", $md)
$wsOverview.Range("B4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C4").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$xlfZhCn = "$uuid.$token.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/407048bd854fcf0d97e6a44f94b8fc8b83b13639/e2e/$md", "", "", $md)
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/407048bd854fcf0d97e6a44f94b8fc8b83b13639/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$xlfZhCn", "", "", $xlfZhCn)
$wsZhCn.Range("E4").Value = "2016-03-20 16:57:12"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/407048bd854fcf0d97e6a44f94b8fc8b83b13639/e2e/$md", "", "", $md)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/407048bd854fcf0d97e6a44f94b8fc8b83b13639/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$xlfZhCn", "", "", $xlfZhCn)
$wsZhCn.Range("H4").Value = "2016-03-20 16:57:59"
$wsZhCn.Range("J4").Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$xlfDeDe = "$uuid.$token.de-de.xlf"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/407048bd854fcf0d97e6a44f94b8fc8b83b13639/e2e/$md", "", "", $md)
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/407048bd854fcf0d97e6a44f94b8fc8b83b13639/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$xlfDeDe", "", "", $xlfDeDe)
$wsDeDe.Range("E4").Value = "2016-03-20 16:57:21"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/407048bd854fcf0d97e6a44f94b8fc8b83b13639/e2e/$md", "", "", $md)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/407048bd854fcf0d97e6a44f94b8fc8b83b13639/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$xlfDeDe", "", "", $xlfDeDe)
$wsDeDe.Range("H4").Value = "2016-03-20 16:58:14"
$wsDeDe.Range("J4").Value = "Include"
